# Populate the "car_price" (J) / "car_comfort_level" (K) columns for rows 2-8
# with the sample data added in this commit, and move the selection to G2
# (clearing the old AI1 scroll position / AM1 selection in the process).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 10).Value = 0   # column J: car_price
    $ws.Cells.Item($row, 11).Value = 1   # column K: car_comfort_level
}

$ws.Range("G2").Select() | Out-Null
